{"js": "// Replace the 25 \"NN\u00f7N=\" division-problem cells in the table with their new values.\n// Matches are resolved against the ORIGINAL document text before any edits are made,\n// so a replacement value that happens to equal another cell's original text\n// (e.g. \"78\u00f75=\" -> \"12\u00f78=\" while a different cell already reads \"12\u00f78=\") cannot\n// be mis-matched against text inserted earlier in this same run.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"58\u00f79=\", \"17\u00f77=\"],\n  [\"28\u00f78=\", \"87\u00f76=\"],\n  [\"59\u00f74=\", \"94\u00f74=\"],\n  [\"47\u00f75=\", \"21\u00f72=\"],\n  [\"80\u00f75=\", \"84\u00f73=\"],\n  [\"75\u00f79=\", \"57\u00f74=\"],\n  [\"41\u00f74=\", \"58\u00f72=\"],\n  [\"39\u00f74=\", \"65\u00f72=\"],\n  [\"89\u00f76=\", \"16\u00f72=\"],\n  [\"54\u00f78=\", \"94\u00f74=\"],\n  [\"79\u00f76=\", \"96\u00f78=\"],\n  [\"44\u00f75=\", \"92\u00f75=\"],\n  [\"42\u00f76=\", \"86\u00f76=\"],\n  [\"37\u00f72=\", \"42\u00f78=\"],\n  [\"78\u00f75=\", \"12\u00f78=\"],\n  [\"85\u00f78=\", \"79\u00f79=\"],\n  [\"63\u00f78=\", \"37\u00f77=\"],\n  [\"49\u00f73=\", \"95\u00f73=\"],\n  [\"71\u00f79=\", \"34\u00f76=\"],\n  [\"89\u00f77=\", \"45\u00f78=\"],\n  [\"86\u00f77=\", \"40\u00f72=\"],\n  [\"31\u00f76=\", \"47\u00f77=\"],\n  [\"55\u00f73=\", \"66\u00f76=\"],\n  [\"31\u00f72=\", \"93\u00f78=\"],\n  [\"12\u00f78=\", \"60\u00f73=\"],\n];\n\n// Phase 1: locate every source range using the pristine document text.\nconst searches = replacements.map(([from, to]) => {\n  const found = body.search(from, { matchCase: true, matchWholeWord: true });\n  found.load(\"items\");\n  return { found, from, to };\n});\nawait context.sync();\n\n// Phase 2: apply the replacements using the ranges captured above.\nfor (const { found, from, to } of searches) {\n  if (found.items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${from}\", found ${found.items.length}`);\n  }\n  found.items[0].insertText(to, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the 25 \"NN\u00f7N=\" division-problem cells in the practice table.\n#\n# $replacements lists the (old, new) text pairs in the same top-to-bottom,\n# left-to-right order the cells appear in the document. One new value\n# (\"12\u00f78=\", produced by row 15: \"78\u00f75=\" -> \"12\u00f78=\") is identical to the old\n# value of a later cell (row 25: \"12\u00f78=\" -> \"60\u00f73=\"). Find.Execute always\n# matches the first occurrence in the document, so the replacements are\n# applied back-to-front (last cell first) - that way a cell is never touched\n# until every cell below it already holds its final text, and \"78\u00f75=\" cannot\n# accidentally re-match the text that row 25 was supposed to receive.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"58\u00f79=\", \"17\u00f77=\")\n    ,@(\"28\u00f78=\", \"87\u00f76=\")\n    ,@(\"59\u00f74=\", \"94\u00f74=\")\n    ,@(\"47\u00f75=\", \"21\u00f72=\")\n    ,@(\"80\u00f75=\", \"84\u00f73=\")\n    ,@(\"75\u00f79=\", \"57\u00f74=\")\n    ,@(\"41\u00f74=\", \"58\u00f72=\")\n    ,@(\"39\u00f74=\", \"65\u00f72=\")\n    ,@(\"89\u00f76=\", \"16\u00f72=\")\n    ,@(\"54\u00f78=\", \"94\u00f74=\")\n    ,@(\"79\u00f76=\", \"96\u00f78=\")\n    ,@(\"44\u00f75=\", \"92\u00f75=\")\n    ,@(\"42\u00f76=\", \"86\u00f76=\")\n    ,@(\"37\u00f72=\", \"42\u00f78=\")\n    ,@(\"78\u00f75=\", \"12\u00f78=\")\n    ,@(\"85\u00f78=\", \"79\u00f79=\")\n    ,@(\"63\u00f78=\", \"37\u00f77=\")\n    ,@(\"49\u00f73=\", \"95\u00f73=\")\n    ,@(\"71\u00f79=\", \"34\u00f76=\")\n    ,@(\"89\u00f77=\", \"45\u00f78=\")\n    ,@(\"86\u00f77=\", \"40\u00f72=\")\n    ,@(\"31\u00f76=\", \"47\u00f77=\")\n    ,@(\"55\u00f73=\", \"66\u00f76=\")\n    ,@(\"31\u00f72=\", \"93\u00f78=\")\n    ,@(\"12\u00f78=\", \"60\u00f73=\")\n)\n\nfor ($i = $replacements.Count - 1; $i -ge 0; $i--) {\n    $oldText = $replacements[$i][0]\n    $newText = $replacements[$i][1]\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $true,      # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Could not find expected text: $oldText\"\n    }\n}\n"}
